$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions) - F-column (want-to-go count) updates only
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 3525
$wsExpo.Range("F5").Value  = 8239
$wsExpo.Range("F7").Value  = 94
$wsExpo.Range("F8").Value  = 2181
$wsExpo.Range("F10").Value = 188
$wsExpo.Range("F12").Value = 1200
$wsExpo.Range("F13").Value = 64
$wsExpo.Range("F15").Value = 27
$wsExpo.Range("F16").Value = 591
$wsExpo.Range("F17").Value = 85
$wsExpo.Range("F18").Value = 6091
$wsExpo.Range("F20").Value = 7313
$wsExpo.Range("F22").Value = 56276
$wsExpo.Range("F23").Value = 4500
$wsExpo.Range("F25").Value = 1042
$wsExpo.Range("F26").Value = 877
$wsExpo.Range("F31").Value = 3732
$wsExpo.Range("F33").Value = 52
$wsExpo.Range("F35").Value = 870
$wsExpo.Range("F36").Value = 1227
$wsExpo.Range("F37").Value = 1216
$wsExpo.Range("F39").Value = 195
$wsExpo.Range("F42").Value = 8
$wsExpo.Range("F43").Value = 769
$wsExpo.Range("F44").Value = 174
$wsExpo.Range("F46").Value = 169
$wsExpo.Range("F47").Value = 6
$wsExpo.Range("F48").Value = 42

# ---------------------------------------------------------------------------
# Sheet "演出" (performances) - F-column updates + a new row inserted at 42
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F12").Value = 115
$wsShow.Range("F14").Value = 41
$wsShow.Range("F16").Value = 7486
$wsShow.Range("F25").Value = 69
$wsShow.Range("F30").Value = 122
$wsShow.Range("F34").Value = 78

# Insert a brand-new row at position 42 (pushes old rows 42-47 down to 43-48)
$wsShow.Rows("42:42").Insert()

# Give the new A42 the same formatting (bold / border / centered) as the
# rest of the index column, then set its value.
$wsShow.Range("A43").Copy()
$wsShow.Range("A42").PasteSpecial(-4122)

$wsShow.Range("A42").Value = 41
$wsShow.Range("B42").NumberFormat = "@"
$wsShow.Range("B42").Value = "2024-11-23"
$wsShow.Range("C42").Value = '上海·w-inds. LIVE TOUR 2024 "Nostalgia" '
$wsShow.Range("D42").Value = "高青西路777号 上海前滩31演艺中心"
$wsShow.Range("E42").Value = "2024.11.23 20:00-11.23 21:30"
$wsShow.Range("F42").Value = 63
$wsShow.Range("G42").Value = 980
$wsShow.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=92863"
$wsShow.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202409/Ol1GyI1x1727235334903.jpeg"

# The row-insert carried the old index numbers down with the rest of each
# row; restore column A so it keeps following the "row number - 1" pattern.
$wsShow.Range("A43").Value = 42
$wsShow.Range("A44").Value = 43
$wsShow.Range("A45").Value = 44
$wsShow.Range("A46").Value = 45
$wsShow.Range("A47").Value = 46
$wsShow.Range("A48").Value = 47

# ---------------------------------------------------------------------------
# Sheet "本地生活" (local life) - F-column updates + a new row appended at 16
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F5").Value  = 1558
$wsLocal.Range("F7").Value  = 663
$wsLocal.Range("F8").Value  = 2342
$wsLocal.Range("F9").Value  = 9339
$wsLocal.Range("F10").Value = 1678
$wsLocal.Range("F12").Value = 89
$wsLocal.Range("F15").Value = 178

$wsLocal.Range("A15").Copy()
$wsLocal.Range("A16").PasteSpecial(-4122)

$wsLocal.Range("A16").Value = 15
$wsLocal.Range("B16").NumberFormat = "@"
$wsLocal.Range("B16").Value = "2024-10-10"
$wsLocal.Range("C16").Value = "上海·「火影忍者疾风传 × animate cafe」"
$wsLocal.Range("D16").Value = "西藏北路198号大悦城北座8楼N809-1 animate cafe上海店"
$wsLocal.Range("E16").Value = "2024.10.10 00:00-11.12 23:59"
$wsLocal.Range("F16").Value = 121
$wsLocal.Range("G16").Value = 30
$wsLocal.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=92883"
$wsLocal.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202409/aQIhaIgt1727249498713.png"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types) - F-column updates + G15 becomes non-sellable
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 3525
$wsAll.Range("F3").Value  = 8239
$wsAll.Range("F4").Value  = 1558
$wsAll.Range("F5").Value  = 2342
$wsAll.Range("F7").Value  = 1678
$wsAll.Range("F9").Value  = 89
$wsAll.Range("F11").Value = 94
$wsAll.Range("G15").Value = "不可售"
$wsAll.Range("F16").Value = 188
$wsAll.Range("F17").Value = 64
$wsAll.Range("F18").Value = 27
$wsAll.Range("F19").Value = 591
$wsAll.Range("F20").Value = 85
$wsAll.Range("F21").Value = 56276
$wsAll.Range("F23").Value = 4500
$wsAll.Range("F27").Value = 115
$wsAll.Range("F29").Value = 3732
$wsAll.Range("F31").Value = 41
$wsAll.Range("F32").Value = 52
$wsAll.Range("F34").Value = 870
$wsAll.Range("F35").Value = 1227
$wsAll.Range("F41").Value = 769
$wsAll.Range("F42").Value = 174
$wsAll.Range("F43").Value = 169
$wsAll.Range("F44").Value = 6
